$d = $word.ActiveDocument

# 1) Remove the "Meta description" paragraph that follows the title heading.
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# 2) Insert a new bold paragraph ("Play Gates of Troy Free Slot Game | Review by
#    SlotJava") right before the final (italic "Create a feature image..." )
#    paragraph. We first open a fresh empty paragraph gap, then replace that
#    gap's contents with the exact run structure we need via InsertXML so
#    that formatting/leading empty run match the target precisely.
$last = $d.Paragraphs($d.Paragraphs.Count)
$gap = $last.Range.Duplicate
$gap.InsertParagraphBefore()

$newPara = $d.Paragraphs($d.Paragraphs.Count - 1)
$newRange = $newPara.Range.Duplicate
$xmlFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Gates of Troy Free Slot Game | Review by SlotJava</w:t></w:r></w:p>'
$newRange.InsertXML($xmlFrag) | Out-Null

# 3) Swap the final paragraph's text (still italic) from the old image prompt
#    to the review blurb.
$finalPara = $d.Paragraphs($d.Paragraphs.Count)
$finalRange = $finalPara.Range.Duplicate
$finalRange.Find.Execute(
    "Create a feature image for " + [char]34 + "Gates of Troy" + [char]34 + " in cartoon style with a happy Maya warrior wearing glasses. The image should showcase the warrior standing in front of the gates of Troy, with soldiers and a wooden horse in the background. The warrior should be holding a large bag of coins and wearing a big smile on their face to emphasize the game's winning potential. The overall tone should be playful and engaging, incorporating bright colors and bold outlines to catch the viewer's attention.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of the Gates of Troy slot game, play for free, and increase your chances of winning big. Optimized for mobile play and features great Greek mythology theme and symbols.",
    2) | Out-Null
